$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting copied from the row above (the bold
# header row). Clear that so the new row looks like a normal data row.
$ws.Rows.Item(2).ClearFormats()

# Restore the date number format on column D (same as the other date cells).
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(2, 1).Value2  = 6
$ws.Cells.Item(2, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(2, 3).Value2  = "Metropolitana"
$ws.Cells.Item(2, 4).Value2  = 45245
$ws.Cells.Item(2, 5).Value2  = 13
$ws.Cells.Item(2, 6).Value2  = "Fruta"
$ws.Cells.Item(2, 7).Value2  = 100108
$ws.Cells.Item(2, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(2, 9).Value2  = 100108007
$ws.Cells.Item(2, 10).Value2 = "Coco"
$ws.Cells.Item(2, 11).Value2 = "Sin especificar"
$ws.Cells.Item(2, 12).Value2 = "Primera"
$ws.Cells.Item(2, 13).Value2 = 50
$ws.Cells.Item(2, 14).Value2 = 30000
$ws.Cells.Item(2, 15).Value2 = 30000
$ws.Cells.Item(2, 16).Value2 = 30000
$ws.Cells.Item(2, 17).Value2 = "$/malla 20 unidades"
$ws.Cells.Item(2, 18).Value2 = "Perú"
$ws.Cells.Item(2, 19).Value2 = 1500
$ws.Cells.Item(2, 20).Value2 = 20
